$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the "Lookup" sheet right after Sheet1 -----------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Lookup"

# Field type options, typed in this order then sorted A-Z (matches final
# alphabetical order seen in the sheet).
$fieldTypes = @("nameOrTitle","description","typeOrCategory","percentageOrRatio","measurement","currency","phoneNumber","emailAddress","orderedOrRanked","binary","locationOrPlaceName","coordinate","countOrAmount","dateAndTime","uniqueIdentifier")
$fieldTypesSorted = $fieldTypes | Sort-Object

$ws2.Range("A1").Value = "Field Type"
$ws2.Range("B1").Value = "Seperator"

for ($i = 0; $i -lt $fieldTypesSorted.Count; $i++) {
    $ws2.Cells.Item($i + 2, 1).Value = $fieldTypesSorted[$i]
}

$ws2.Range("B2").Value = "Yes"
$ws2.Range("B3").Value = "No"

$ws2.Columns.Item(1).ColumnWidth = 19.59
$ws2.Columns.Item(2).ColumnWidth = 10.92

# --- Turn the two lookup ranges into Excel Tables --------------------------
$tbl1 = $ws2.ListObjects.Add(1, $ws2.Range("A1:A16"), $null, 1)
$tbl1.Name = "Table1"

$tbl2 = $ws2.ListObjects.Add(1, $ws2.Range("B1:B3"), $null, 1)
$tbl2.Name = "Table2"

# --- Defined names used by the data validation drop-downs ------------------
$wb.Names.Add("FieldType", "=Lookup!`$A`$2:`$A`$16")
$wb.Names.Add("Seperator", "=Lookup!`$B`$2:`$B`$3")

# --- Data on Sheet1: fill in the thousandsSeparator column ----------------
$ws1.Range("F2").Value = "No"
$ws1.Range("F3").Value = "No"
$ws1.Range("F4").Value = "Yes"

# --- Data validation drop-downs on Sheet1 -----------------------------------
$dv1 = $ws1.Range("D1:D1048576").Validation
$dv1.Add(3, 1, 1, "=FieldType")
$dv1.InputTitle = "Instructions"
$dv1.InputMessage = "Select a data type from the list."
$dv1.ShowInput = $true
$dv1.ShowError = $true

$dv2 = $ws1.Range("F1:F1048576").Validation
$dv2.Add(3, 1, 1, "=Seperator")
$dv2.InputTitle = "Instructions"
$dv2.InputMessage = "Select a value from the list."
$dv2.ShowInput = $true
$dv2.ShowError = $true

# --- View tweaks on Sheet1 (matches author re-selecting a different cell) --
$ws1.Range("F4").Select()
$excel.ActiveWindow.ScrollColumn = 2

$ws1.Activate()
